$d = $word.ActiveDocument

# Helper: replace the first occurrence of $find inside the given Range with
# $replace, without letting the engine "smear" a neighbouring run's
# formatting (or displace nearby range markers such as comment anchors) onto
# the new text - both of which have been observed when using a plain
# Find.Execute(..., Replace:=wdReplaceAll) on text that sits right next to a
# run/hyperlink/comment-range boundary. We locate the match, insert the new
# text collapsed at the match's own end (so it inherits only the format that
# is already active there, inside the original run, and any boundary marker
# immediately before the match keeps pointing at the same spot), then delete
# the old text that is now right before it.
function Replace-InRange {
    param($scope, $find, $replace)
    $r = $scope.Duplicate
    $found = $r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $s = $r.Start
        $e = $r.End
        $ins = $d.Range($e, $e)
        $ins.InsertBefore($replace)
        $old = $d.Range($s, $e)
        $old.Delete()
    }
    return $found
}

# --- Hyperlink display texts (use TextToDisplay so the <w:hyperlink> wrapper
#     and the run's own rPr - colour/underline - survive intact) -----------
foreach ($h in $d.Hyperlinks) {
    if ($h.TextToDisplay -eq "English") {
        $h.TextToDisplay = "ภาษาอังกฤษ"
    } elseif ($h.TextToDisplay -eq "live chat") {
        $h.TextToDisplay = "แชทสด"
    }
}

# --- Plain "English" label paragraph (not a hyperlink) ---------------------
Replace-InRange $d.Content "English" "ภาษาอังกฤษ" | Out-Null

# --- language list line, right after the "English" hyperlink run ----------
Replace-InRange $d.Content " / Portuguese / French / Thai / Vietnamese / Spanish" " / ภาษาโปรตุเกส / ภาษาฝรั่งเศส /ภาษาไทย / ภาษาเวียดนาม / ภาษาสเปน" | Out-Null

# --- "Brief" ----------------------------------------------------------------
Replace-InRange $d.Content "Brief" "บทย่อ" | Out-Null

# --- Brief description -------------------------------------------------------
Replace-InRange $d.Content "An email sent upon verification to partners in the target country who have sent the correct documents. It will be sent via customer.io" "อีเมลที่ส่งเมื่อได้รับการยืนยันไปยังพันธมิตรในประเทศเป้าหมายที่ได้ส่งเอกสารที่ถูกต้องแล้ว โดยมันจะถูกส่งผ่านทาง customer.io" | Out-Null

# --- "Target audience" ------------------------------------------------------
Replace-InRange $d.Content "Target audience" "กลุ่มเป้าหมาย" | Out-Null

# --- Target audience description --------------------------------------------
Replace-InRange $d.Content "Invited partners who didn’t submit their documents on time" "พันธมิตรที่ได้รับเชิญซึ่งไม่ได้ส่งเอกสารตรงเวลา" | Out-Null

# --- "Subject line" ----------------------------------------------------------
Replace-InRange $d.Content "Subject line" "หัวเรื่อง" | Out-Null

# --- subject line tail --------------------------------------------------------
Replace-InRange $d.Content " — one step closer!" " — อีกหนึ่งก้าวใกล้แล้ว!" | Out-Null

# --- headline -------------------------------------------------------------
Replace-InRange $d.Content "Your documents have been verified!" "เอกสารของคุณได้รับการตรวจสอบยืนยันแล้ว!" | Out-Null

# --- "Hi " -------------------------------------------------------------------
Replace-InRange $d.Content "Hi " "สวัสดี " | Out-Null

# --- ", " right after [PARTNER NAME] becomes just " " (paragraph 16) --------
Replace-InRange $d.Paragraphs.Item(16).Range ", " " " | Out-Null

# --- "We’ve reviewed..." -------------------------------------------------------
Replace-InRange $d.Content "We’ve reviewed the documents you’ve sent us for the " "เราได้ตรวจสอบเอกสารที่คุณส่งมาให้เราสำหรับงาน " | Out-Null

# --- "and all of them have been verified!" --------------------------------------
Replace-InRange $d.Content " and all of them have been verified! " " และเอกสารทั้งหมดได้รับการตรวจสอบยืนยันเรียบร้อยแล้ว! " | Out-Null

# --- "We’ll be sending..." ---------------------------------------------------
Replace-InRange $d.Content "We’ll be sending out more details about the event soon, including the agenda and travel arrangements, so make sure to check your inbox regularly." "เราจะส่งรายละเอียดเพิ่มเติมเกี่ยวกับกิจกรรมไปให้คุณในเร็วๆ นี้ รวมถึงกำหนดการและการเตรียมเรื่องการเดินทาง ดังนั้นโปรดตรวจดูกล่องข้อความอีเมล์ของคุณอย่างสม่ำเสมอ" | Out-Null

# --- "If you have any questions, please contact us via " -----------------------
Replace-InRange $d.Content "If you have any questions, please contact us via " "หากคุณมีคำถามใดๆ กรุณาติดต่อเราผ่านทาง " | Out-Null

# --- " or " between the live-chat and WhatsApp hyperlinks (paragraph 20) ----
Replace-InRange $d.Paragraphs.Item(20).Range " or " " หรือทาง " | Out-Null

# --- ". " right after the WhatsApp hyperlink -> " " (still paragraph 20) ---
Replace-InRange $d.Paragraphs.Item(20).Range ". " " " | Out-Null

# --- "If you have any questions, please contact your country manager, " ----
Replace-InRange $d.Content "If you have any questions, please contact your country manager, " "หากคุณมีคำถามใดๆ โปรดติดต่อผู้จัดการประจำประเทศของคุณซึ่งได้แก่ " | Out-Null

# --- ", at " -------------------------------------------------------------------
Replace-InRange $d.Content ", at " " ที่ " | Out-Null

# --- " or " between [EMAIL ADDRESS] and [WHATSAPP NO] (paragraph 21) -------
Replace-InRange $d.Paragraphs.Item(21).Range " or " " หรือ " | Out-Null

# --- " (WhatsApp). " ------------------------------------------------------------
Replace-InRange $d.Content " (WhatsApp). " " (WhatsApp) " | Out-Null

# --- comment text "choose either one" ---------------------------------------
foreach ($c in $d.Comments) {
    Replace-InRange $c.Range "choose either one" "เลือกอย่างใดอย่างหนึ่ง" | Out-Null
}
